$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259, shifting existing rows 259:289 down to 260:290.
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly data point.
$ws.Range("A259").Value = 10
$ws.Range("B259").Value = "Vega Modelo de Temuco"
$ws.Range("C259").Value = "La Araucanía"
$ws.Range("D259").Value = 45124
$ws.Range("E259").Value = 9
$ws.Range("F259").Value = 100112012
$ws.Range("G259").Value = "Espinaca"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 80
$ws.Range("K259").Value = 8000
$ws.Range("L259").Value = 8000
$ws.Range("M259").Value = 8000
$ws.Range("N259").Value = "`$/docena de atados"
$ws.Range("O259").Value = "Región de La Araucanía"
$ws.Range("P259").Value = 2667
$ws.Range("Q259").Value = 3
$ws.Range("R259").Value = "Hortaliza"
